$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) "Actualiza base de datos EC": the period column (E16:E24, "Periodo Mora")
#    used to list periods in descending order (1703 .. 1607). The refreshed
#    export lists them ascending (1607 .. 1703).
# ---------------------------------------------------------------------------
$periods = @("1607", "1608", "1609", "1610", "1611", "1612", "1701", "1702", "1703")
for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Range("E$row").Value = $periods[$i]
}

# ---------------------------------------------------------------------------
# 2) Column widths were re-measured by the newer Excel build that produced
#    the refreshed workbook (the data in column C, "N deg Doc Trabajador",
#    got narrower; the rest shifted slightly too). Reapply the new best-fit
#    widths captured from the updated file.
# ---------------------------------------------------------------------------
$colWidths = @{
    2  = 16.072916666666668
    3  = 7.983072916666667
    4  = 18.709635416666668
    5  = 11.893229166666666
    6  = 8.619791666666666
    7  = 12.619791666666666
    8  = 17.072916666666668
    9  = 15.983072916666666
    10 = 13.346354166666666
}
foreach ($col in $colWidths.Keys) {
    $ws.Columns.Item($col).ColumnWidth = $colWidths[$col]
}

# ---------------------------------------------------------------------------
# 3) The company logo picture shifted left (anchor moved, same size) once the
#    new "Parte 1" layout was applied.
# ---------------------------------------------------------------------------
$shp = $ws.Shapes.Item(1)
$shp.Width = 76.81889763779527
$shp.Height = 48.18897637795275
$shp.Left = 53.59055118110236
$shp.Top = 19.405511811023622
